$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.453.45"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3
$ws.Range("D3").Value = "1.905.76"
$ws.Range("E3").Value = "  +2.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.86"
$ws.Range("E5").Value = "  +3.56%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.639"
$ws.Range("E6").Value = "  +2.71%  "

# Row 7
$ws.Range("E7").Value = "  +0.19%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.88"
$ws.Range("E8").Value = "  -1.43%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.341"
$ws.Range("E9").Value = "  +3.89%  "

# Row 10
$ws.Range("E10").Value = "  +1.30%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0998"
$ws.Range("E11").Value = "  +0.97%  "

# Row 12
$ws.Range("D12").Value = "2.182.78"
$ws.Range("E12").Value = "  +2.83%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.34"
$ws.Range("E13").Value = "  +7.98%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.697"
$ws.Range("E14").Value = "  +2.69%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.881.12"
$ws.Range("E15").Value = "  +1.17%  "

# Row 16
$ws.Range("E16").Value = "  +2.55%  "

# Row 17
$ws.Range("D17").Value = "35.498.08"
$ws.Range("E17").Value = "  +1.26%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.89"

# Row 19
$ws.Range("E19").Value = "  +3.69%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.21"
$ws.Range("E20").Value = "  +0.82%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.69"
$ws.Range("E21").Value = "  +4.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.85"
$ws.Range("E22").Value = "  +2.14%  "

# Row 23
$ws.Range("E23").Value = "  +0.19%  "

# Row 24
$ws.Range("E24").Value = "  +0.26%  "

# Row 25
$ws.Range("E25").Value = "  +0.46%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("E26").Value = "  +12.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.56"
$ws.Range("E27").Value = "  +8.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.97"
$ws.Range("E28").Value = "  +1.81%  "

# Row 29
$ws.Range("E29").Value = "  +1.09%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.977"
$ws.Range("E30").Value = "  +26.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0570"
$ws.Range("E31").Value = "  +2.34%  "

# Row 32
$ws.Range("E32").Value = "  +2.98%  "

# Row 33
$ws.Range("E33").Value = "  +0.35%  "

# Row 34
$ws.Range("E34").Value = "  +4.47%  "

# Row 35
$ws.Range("E35").Value = "  +8.47%  "

# Row 36
$ws.Range("E36").Value = "  +0.40%  "

# Row 37
$ws.Range("E37").Value = "  +7.74%  "

# Row 38
$ws.Range("E38").Value = "  +2.95%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0205"
$ws.Range("E39").Value = "  +1.30%  "

# Row 40
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "91.90"
$ws.Range("E40").Value = "  +0.42%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0633"
$ws.Range("E41").Value = "  +14.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.60"
$ws.Range("E42").Value = "  +4.25%  "

# Row 43
$ws.Range("D43").Value = "1.349.32"
$ws.Range("E43").Value = "  +0.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "50.32"
$ws.Range("E44").Value = "  +45.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  +2.25%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.88"
$ws.Range("E46").Value = "  +0.96%  "

# Row 47
$ws.Range("E47").Value = "  +0.14%  "

# Row 48
$ws.Range("E48").Value = "  +0.18%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.64"
$ws.Range("E49").Value = "  +4.76%  "

# Row 50
$ws.Range("D50").Value = "2.093.55"
$ws.Range("E50").Value = "  +2.54%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0692"
$ws.Range("E51").Value = "  +1.59%  "
